# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values. Values in column D that look like
# plain numbers are forced to Text format first so Excel keeps them
# as literal strings (matching the source inlineStr cells) instead of
# silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.348.56"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.841.86"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.93"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6299"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07728"

$ws.Range("D12").Value = "1.843.86"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.967"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6756"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001027"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.58"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.245"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").Value = "29.380.91"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.01"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.365"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.97"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.489"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1348"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06906"
$ws.Range("E28").Value = "  +8.16%  "

$ws.Range("E29").Value = "  +4.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.046"
$ws.Range("E31").Value = "  -1.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6975"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.586"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01841"
$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.818"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").Value = "1.232.72"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.804"
$ws.Range("E40").Value = "  +3.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9354"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "1.992.20"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.93"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.23"
$ws.Range("E45").Value = "  -1.64%  "

$ws.Range("E46").Value = "  +4.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.020"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.703"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.956"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1139"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3906"
$ws.Range("E51").Value = "  -0.94%  "
